$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Text)
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.340.92"
Set-TextValue "E2" "  -2.26%  "
Set-TextValue "D3" "1.793.81"
Set-TextValue "E3" "  -2.01%  "
Set-TextValue "D4" "1.005"
Set-TextValue "E4" "  -0.30%  "
Set-TextValue "D5" "1.004"
Set-TextValue "E5" "  -0.32%  "
Set-TextValue "D6" "306.62"
Set-TextValue "E6" "  -1.44%  "
Set-TextValue "D7" "0.4509"
Set-TextValue "E7" "  -1.50%  "
Set-TextValue "D8" "0.3596"
Set-TextValue "E8" "  -2.94%  "
Set-TextValue "D9" "45.90"
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "0.07073"
Set-TextValue "E10" "  -1.44%  "
Set-TextValue "D11" "0.8827"
Set-TextValue "E11" "  +0.75%  "
Set-TextValue "D12" "0.07750"
Set-TextValue "E12" "  -0.22%  "
Set-TextValue "D13" "19.45"
Set-TextValue "E13" "  -0.97%  "
Set-TextValue "D14" "1.812.53"
Set-TextValue "E14" "  -2.39%  "
Set-TextValue "D15" "5.275"
Set-TextValue "E15" "  -1.15%  "
Set-TextValue "D16" "6.312"
Set-TextValue "E16" "  -1.34%  "
Set-TextValue "D17" "84.81"
Set-TextValue "E17" "  -2.88%  "
Set-TextValue "E18" "  -0.31%  "
Set-TextValue "D19" "0.000008504"
Set-TextValue "E19" "  -2.49%  "
Set-TextValue "D20" "1.004"
Set-TextValue "E20" "  -0.25%  "
Set-TextValue "E21" "  -1.89%  "
Set-TextValue "D22" "26.368.61"
Set-TextValue "E22" "  -2.27%  "
Set-TextValue "B24" "WrappedliquidstakedEther2.0"
Set-TextValue "C24" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D24" "2.060.46"
Set-TextValue "E24" "  -0.29%  "
Set-TextValue "B25" "Cosmos"
Set-TextValue "C25" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D25" "10.52"
Set-TextValue "E25" "  +0.81%  "
Set-TextValue "D26" "1.968"
Set-TextValue "E26" "  -1.98%  "
Set-TextValue "D27" "150.97"
Set-TextValue "E27" "  -0.26%  "
Set-TextValue "D28" "17.80"
Set-TextValue "E28" "  -2.25%  "
Set-TextValue "D29" "2.006"
Set-TextValue "E29" "  +2.08%  "
Set-TextValue "D30" "111.87"
Set-TextValue "E30" "  -1.84%  "
Set-TextValue "D31" "4.867"
Set-TextValue "E31" "  -1.35%  "
Set-TextValue "D32" "0.08663"
Set-TextValue "E32" "  -1.69%  "
Set-TextValue "D33" "3.068"
Set-TextValue "E33" "  +2.96%  "
Set-TextValue "E34" "  -1.21%  "
Set-TextValue "B35" "RenderToken"
Set-TextValue "C35" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D35" "2.715"
Set-TextValue "E35" "  +5.78%  "
Set-TextValue "B36" "ImmutableX"
Set-TextValue "C36" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D36" "0.7226"
Set-TextValue "E36" "  -3.73%  "
Set-TextValue "D37" "1.104"
Set-TextValue "E37" "  -2.87%  "
Set-TextValue "E38" "  -0.20%  "
Set-TextValue "D39" "1.065"
Set-TextValue "E39" "  -2.12%  "
Set-TextValue "D40" "0.01929"
Set-TextValue "E40" "  -1.09%  "
Set-TextValue "D41" "0.05085"
Set-TextValue "E41" "  -1.25%  "
Set-TextValue "D42" "2.863"
Set-TextValue "E42" "  -1.21%  "
Set-TextValue "D43" "0.5034"
Set-TextValue "E43" "  +1.30%  "
Set-TextValue "D44" "6.841"
Set-TextValue "E44" "  -1.70%  "
Set-TextValue "E45" "  -5.55%  "
Set-TextValue "D46" "7.987"
Set-TextValue "E46" "  -4.12%  "
Set-TextValue "D47" "1.004"
Set-TextValue "E47" "  -0.36%  "
Set-TextValue "D48" "0.4613"
Set-TextValue "E48" "  -1.60%  "
Set-TextValue "B49" "Quant"
Set-TextValue "C49" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D49" "101.04"
Set-TextValue "E49" "  -1.29%  "
Set-TextValue "B50" "EnergySwap"
Set-TextValue "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "9.812"
Set-TextValue "E50" "  -3.10%  "
Set-TextValue "E51" "  -2.36%  "
